$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared by all rows in this block (Agrícola del Norte S.A. de Arica - Pepino dulce)
$colA = 1
$colB = "Agrícola del Norte S.A. de Arica"
$colC = "Arica y Parinacota"
$colE = 15
$colF = 100112043
$colG = "Pepino dulce"
$colH = "Cultivar XV región"
$colO = "Región de Arica y Parinacota"
$colR = "Hortaliza"

# New data for rows 32..38 (columns: D,I,J,K,L,M,N,P,Q)
$rows = @(
    @{ Row = 32; D = 45212; I = "Segunda"; J = 100; K = 19000; L = 20000; M = 19500; N = "`$/bandeja 18 kilos"; P = 1083; Q = 18 },
    @{ Row = 33; D = 45212; I = "Tercera"; J = 120; K = 16000; L = 17000; M = 16500; N = "`$/bandeja 18 kilos"; P = 917;  Q = 18 },
    @{ Row = 34; D = 45212; I = "Tercera"; J = 160; K = 6000;  L = 7000;  M = 6500;  N = "`$/caja 10 kilos";    P = 650;  Q = 10 },
    @{ Row = 35; D = 44526; I = "Primera"; J = 100; K = 5000;  L = 5500;  M = 5250;  N = "`$/caja 10 kilos";    P = 525;  Q = 10 },
    @{ Row = 36; D = 44526; I = "Segunda"; J = 100; K = 4000;  L = 4500;  M = 4250;  N = "`$/caja 10 kilos";    P = 425;  Q = 10 },
    @{ Row = 37; D = 44526; I = "Tercera"; J = 120; K = 3000;  L = 3500;  M = 3250;  N = "`$/caja 10 kilos";    P = 325;  Q = 10 },
    @{ Row = 38; D = 44757; I = "Primera"; J = 150; K = 6000;  L = 6500;  M = 6250;  N = "`$/caja 10 kilos";    P = 625;  Q = 10 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $colO
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $colR
}
